# Set line spacing to single (1.0) for every paragraph in the document.
# This maps to <w:spacing w:line="240" w:lineRule="auto"/> in the
# paragraph's <w:pPr> (creating the <w:pPr> if it doesn't already exist).
$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    $p.LineSpacingRule = 0   # wdLineSpaceSingle
}
